# Coinranking symbol-list refresh (Mon Jan  2 05:51:40 UTC 2023 GitHub Actions run).
# Columns D (Price) and E (Volume 1h) hold plain numeric-/percent-looking text
# (t="inlineStr" in the source, i.e. literal strings, not numbers). Excel's COM
# layer auto-coerces a bare numeric/percent-looking string into a real number,
# so each such write is prefixed with a leading apostrophe (forces text) and the
# cell style is put back to "Normal" right after so no stray NumberFormat /
# quotePrefix styling sticks around on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "243.83"
Set-TextValue "E2" "-0.06%"

# Row 3 - OKB
Set-TextValue "D3" "30.00"
Set-TextValue "E3" "13.86%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.141"
Set-TextValue "E4" "-0.17%"

# Row 5 - Cronos
Set-TextValue "D5" "0.05672"
Set-TextValue "E5" "1.39%"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.523"
Set-TextValue "E6" "0.79%"

# Row 7 - MXToken
Set-TextValue "D7" "0.8434"
Set-TextValue "E7" "2.84%"

# Row 8 - FTXToken
Set-TextValue "D8" "0.8576"
Set-TextValue "E8" "2.66%"

# Row 9 - WazirX
Set-TextValue "E9" "0.57%"

# Row 10 - MandalaExchangeToken
Set-TextValue "D10" "0.06910"
Set-TextValue "E10" "-1.19%"

# Row 11 - BitrueCoin
Set-TextValue "D11" "0.02888"
Set-TextValue "E11" "0.03%"

# Row 12 - BitMartToken
Set-TextValue "D12" "0.09387"
Set-TextValue "E12" "0.01%"

# Row 13 - BitForexToken
Set-TextValue "D13" "0.001526"
Set-TextValue "E13" "0.80%"

# Row 14 - CoinExToken
Set-TextValue "D14" "0.04164"
Set-TextValue "E14" "-10.10%"

# Rows 15-19 shuffled: new ranking order is One, TigerCash, LEO, GateToken, BTSEToken
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006013"
Set-TextValue "E15" "0.38%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006096"
Set-TextValue "E16" "-1.63%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.508"
Set-TextValue "E17" "-4.06%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D18" "3.022"
Set-TextValue "E18" "-0.33%"

$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D19" "2.133"
Set-TextValue "E19" "-2.28%"

# Row 20 - BitpandaEcosystemToken
Set-TextValue "E20" "1.20%"

# Row 21 - LiechtensteinCryptoassetsExchange
Set-TextValue "D21" "0.03270"
Set-TextValue "E21" "5.35%"

# Row 22 - ProBitToken
Set-TextValue "E22" "0.24%"

# Row 23 - MCDex
Set-TextValue "D23" "3.619"
Set-TextValue "E23" "-3.18%"

# Row 25 - BitKan
Set-TextValue "D25" "0.001209"
Set-TextValue "E25" "-3.17%"

# Row 26 - HotbitToken
Set-TextValue "D26" "0.004441"
Set-TextValue "E26" "-1.13%"

# Row 27 - NitroEx
Set-TextValue "E27" "22.82%"

# Row 28 - UpBots
Set-TextValue "D28" "0.0001396"
Set-TextValue "E28" "0.22%"

# Row 40 - IDEX
Set-TextValue "D40" "0.03710"
Set-TextValue "E40" "1.84%"

# Rows 41-42 swapped: new order is BKEXToken, KickToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1057"
Set-TextValue "E41" "-23.06%"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.003416"
Set-TextValue "E42" "-44.45%"

# Row 43 - CEJI
Set-TextValue "D43" "0.002287"
Set-TextValue "E43" "-13.03%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.009797"
Set-TextValue "E44" "8.85%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005090"
Set-TextValue "E45" "-4.59%"

# Row 46 - Kangarootoken
Set-TextValue "E46" "-0.12%"

# Row 47 - CoinbaseStockToken
Set-TextValue "D47" "0.09989"
Set-TextValue "E47" "-30.63%"

# Row 48 - BOLO
Set-TextValue "D48" "0.002809"
Set-TextValue "E48" "21.65%"

# Row 49 - CryptobidCoin
Set-TextValue "E49" "-0.12%"

# Row 50 - SpecialPowerGold
Set-TextValue "E50" "-0.12%"
